# MILESTONE: distritos split by concept
# Column A ("largo") is renamed to "lista" and every numeric district code in
# column A (rows 2-74) is replaced by the corresponding list/"lista" name,
# matching the bold/bordered/centered header style already used by row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "lista"

# Grab the header cell's format (bold font, thin border, centered/top aligned)
# so it can be stamped onto each rewritten data cell without inventing a new style.
$ws.Range("A1").Copy()

$ws.Range("A2").Value = "A PULSO, POR EL BUEN VIVIR "
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").Value = "APRUEBO DIGNIDAD"
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").Value = "ARICA SIEMPRE ARICA "
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").Value = "ASAMBLEA CONSTITUYENTE ATACAMA"
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").Value = "ASAMBLEA POPULAR CONSTITUYENTE "
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").Value = "ASAMBLEA POPULAR POR LA DIGNIDAD "
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").Value = "AUTONOMIA SOCIAL Y SINDICAL TARAPACA "
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").Value = "BIOBIO SIN PARTIDOS "
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").Value = "CABILDO AUTOCONVOCADO "
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").Value = "CANDIDATURA INDEPENDIENTE"
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").Value = "CHILE INDEPENDIENTE "
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A13").Value = "CIUDADANOS CRISTIANOS"
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").Value = "COMUNIDAD INDEPENDIENTE DE MAULE "
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A15").Value = "COMUNIDAD INDEPENDIENTE VENSEREMOS "
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A16").Value = "COORDINADORA SOCIAL DE MAGALLANES "
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").Value = "CORRIENTES INDEPENDIENTES "
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").Value = "DECISION CIUDADANA "
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").Value = "ELIGE LA LISTA DEL PUEBLO "
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A20").Value = "ENERGIA INDEPENDIENTE "
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").Value = "FUERZA SOCIAL DE ÑUBLE, LA LISTA DEL PUEBLO "
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").Value = "INDEPENDIENTES COMO TU "
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23").Value = "INDEPENDIENTES CON CHILE "
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").Value = "INDEPENDIENTES DE TARAPACA "
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").Value = "INDEPENDIENTES DE ÑUBLE POR LA NUEVA CONSTITUCION "
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A26").Value = "INDEPENDIENTES DEL APRUEBO REGION COQUIMBO"
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A27").Value = "INDEPENDIENTES DEL BIOBIO POR UNA NUEVA CONSTITUCION "
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A28").Value = "INDEPENDIENTES DEL NORTE GRANDE POR UNA NUEVA CONSTITUCION "
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A29").Value = "INDEPENDIENTES DISTRITO 6 + LISTA DEL PUEBLO "
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A30").Value = "INDEPENDIENTES NUEVA CONSTITUCION "
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A31").Value = "INDEPENDIENTES POR LA NUEVA CONSTITUCION "
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A32").Value = "INDEPENDIENTES POR LA REGION DE COQUIMBO "
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A33").Value = "INDEPENDIENTES POR UNA NUEVA CONSTITUCION "
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A34").Value = "INDEPENDIENTES SIN PADRINOS "
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A35").Value = "INDEPENDIENTES Y MOVIMIENTOS SOCIALES DEL APRUEBO "
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A36").Value = "INSULARES E INDEPENDIENTES "
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A37").Value = "LA LISTA DEL PUEBLO"
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A38").Value = "LA LISTA DEL PUEBLO "
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A39").Value = "LA LISTA DEL PUEBLO 100% INDEPENDIENTES "
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A40").Value = "LA LISTA DEL PUEBLO DISTRITO 12 "
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A41").Value = "LA LISTA DEL PUEBLO DISTRITO 14 "
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A42").Value = "LA LISTA DEL PUEBLO DISTRITO 9 "
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A43").Value = "LA LISTA DEL PUEBLO MAULE SUR "
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("A44").Value = "LISTA DEL APRUEBO"
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A45").Value = "LISTA DEL PUEBLO TRANSFORMANDO DESDE EL WILLI "
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A46").Value = "LISTA DEL PUEBLO – MOVIMIENTO TERRITORIAL CONSTITUYENTE "
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A47").Value = "LISTA DEL PUEBLO-RIOS INDEPENDIENTES "
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A48").Value = "LISTA INDEPENDIENTE PARTO SOCIAL "
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A49").Value = "LISTA POR LA JUSTICIA SOCIAL "
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A50").Value = "LISTA SOCIAL PODER CONSTITUYENTE A TODA COSTA "
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A51").Value = "MAGALLANICOS NO NEUTRALES "
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("A52").Value = "MOVIMIENTO INDEPENDIENTES DEL NORTE "
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("A53").Value = "MOVIMIENTO SOCIAL CONSTITUYENTE / LA LISTA DEL PUEBLO "
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A54").Value = "MOVIMIENTO SOCIAL LA LISTA DEL PUEBLO "
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("A55").Value = "MOVIMIENTOS SOCIALES : UNIDAD DE INDEPENDIENTES "
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("A56").Value = "MOVIMIENTOS SOCIALES AUTONOMOS "
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("A57").Value = "MOVIMIENTOS SOCIALES INDEPENDIENTES "
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A58").Value = "MOVIMIENTOS SOCIALES PLURINACIONALES E INDEPENDIENTES "
$ws.Range("A58").PasteSpecial(-4122)
$ws.Range("A59").Value = "NOBLES HIJXOS DE TARAPACA "
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("A60").Value = "NUESTRAS VOCES "
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("A61").Value = "ORGANIZACIONES SOCIALES Y TERRITORIALES DEL WALLMAPU "
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("A62").Value = "PARTIDO DE TRABAJADORES REVOLUCIONARIOS"
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A63").Value = "PARTIDO ECOLOGISTA VERDE"
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A64").Value = "PARTIDO HUMANISTA"
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A65").Value = "PARTIDO UNION PATRIOTICA"
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A66").Value = "PATAGONIA SOMOS TODOS "
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("A67").Value = "PUEBLO UNIDO TARAPACA "
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("A68").Value = "REGIONALISMO CIUDADANO INDEPENDIENTE "
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("A69").Value = "REPUBLICA DE LOS INDEPENDIENTES DE MAGALLANES "
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("A70").Value = "SEXTA UNIDA "
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("A71").Value = "SOBERANIA CIUDADANA "
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("A72").Value = "SUMATE AHORA "
$ws.Range("A72").PasteSpecial(-4122)
$ws.Range("A73").Value = "VAMOS POR CHILE"
$ws.Range("A73").PasteSpecial(-4122)
$ws.Range("A74").Value = "VOCES CONSTITUYENTES "
$ws.Range("A74").PasteSpecial(-4122)
